$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns D (Price) and E (Volume(1h)) contain numeric-looking text that must
# stay stored as TEXT (matches the original inlineStr cells). Force the text
# number format before assigning, then restore the base style so no stray
# formatting is introduced.
function Set-TextValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

Set-TextValue 'D2' '304.51'
Set-TextValue 'E2' '3.89%'
Set-TextValue 'D3' '32.33'
Set-TextValue 'E3' '6.08%'
Set-TextValue 'D4' '5.309'
Set-TextValue 'E4' '3.00%'
Set-TextValue 'D5' '0.07605'
Set-TextValue 'E5' '6.67%'
Set-TextValue 'D6' '7.885'
Set-TextValue 'E6' '4.72%'
Set-TextValue 'D7' '3.875'
Set-TextValue 'E7' '7.19%'
Set-TextValue 'D8' '1.708'
Set-TextValue 'E8' '22.05%'
Set-TextValue 'D9' '0.9287'
Set-TextValue 'E9' '1.21%'
Set-TextValue 'D10' '0.1701'
Set-TextValue 'E10' '3.89%'
Set-TextValue 'D11' '0.08007'
Set-TextValue 'E11' '1.82%'
Set-TextValue 'D12' '0.08040'
Set-TextValue 'E12' '3.70%'
Set-TextValue 'D13' '0.03062'
Set-TextValue 'E13' '3.98%'
Set-TextValue 'E14' '10.39%'
Set-TextValue 'D15' '0.001497'
Set-TextValue 'E15' '-5.12%'
Set-TextValue 'D16' '0.04593'
Set-TextValue 'E16' '1.20%'
Set-TextValue 'D17' '0.006326'
Set-TextValue 'E17' '0.72%'
Set-TextValue 'E18' '-1.13%'
Set-TextValue 'D19' '2.238'
Set-TextValue 'E19' '-0.20%'
Set-TextValue 'E20' '1.51%'
Set-TextValue 'D21' '0.1344'
Set-TextValue 'E21' '-1.46%'
Set-TextValue 'D22' '4.551'
Set-TextValue 'E22' '9.66%'
Set-TextValue 'D23' '0.1616'
Set-TextValue 'E23' '1.66%'
Set-TextValue 'D24' '0.001217'
Set-TextValue 'E24' '0.50%'
Set-TextValue 'D25' '0.004495'
Set-TextValue 'E25' '5.99%'
Set-TextValue 'D26' '0.0001397'
Set-TextValue 'E26' '19.47%'
Set-TextValue 'D27' '0.0001698'
Set-TextValue 'E27' '0.57%'
Set-TextValue 'D39' '0.01737'
Set-TextValue 'E39' '2,540.99%'
Set-TextValue 'D40' '0.04543'
Set-TextValue 'E40' '2.90%'
Set-TextValue 'D41' '0.006957'
Set-TextValue 'E41' '-1.22%'
Set-TextValue 'D42' '0.1363'
Set-TextValue 'E42' '6.92%'
$ws.Range('B43').Value = 'LocalTraders'
$ws.Range('C43').Value = 'https://coinranking.com/coin/E6DwMU2zXb+localtraders-lct'
Set-TextValue 'D43' '0.01392'
Set-TextValue 'E43' '5.04%'
$ws.Range('B44').Value = 'CEJI'
$ws.Range('C44').Value = 'https://coinranking.com/coin/SbKjCVJCh+ceji-ceji'
Set-TextValue 'D44' '0.002065'
Set-TextValue 'E44' '-6.48%'
Set-TextValue 'D45' '0.00006156'
Set-TextValue 'E45' '4.96%'
Set-TextValue 'D46' '0.7194'
Set-TextValue 'E46' '-58.31%'
Set-TextValue 'D47' '0.01297'
Set-TextValue 'E47' '-0.17%'
